# Backup QR Scanner data - 4/5/2025, 9:52:30 PM
# Adds a new worksheet "Mazinjsbdb" (a fresh QR-scan log) at the end of the
# workbook, mirroring the structure used by the other scan-log sheets.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab order (matching the appended <sheet> entry in the diff).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Mazinjsbdb"

# Force the columns that must stay textual (IDs, dates, times) to a text
# number format before writing, so Excel doesn't silently convert them to
# numbers/dates.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("D1:E2").NumberFormat = "@"

# Header row.
$ws.Range("A1").Value = "Number"
$ws.Range("B1").Value = "Student ID"
$ws.Range("C1").Value = "Location"
$ws.Range("D1").Value = "Log Date"
$ws.Range("E1").Value = "Log Time"

# Data row.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "231249"
$ws.Range("C2").Value = "Mazinjsbdb"
$ws.Range("D2").Value = "2025-04-05"
$ws.Range("E2").Value = "21:52:27"
